$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.442.03"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.58"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.86"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.96"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3389"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07567"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.145"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.020"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.952"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.85"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001123"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.02"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06755"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.305"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.29"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.16"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.433.57"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.333"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.679"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.50"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.034"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.748.23"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.057"
$ws.Range("E32").Value = "  +7.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.168"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.986"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.831"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08371"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02478"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.357"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2306"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.464"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.32"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6221"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.811"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "130.46"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5809"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.065"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.220"
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07316"
$ws.Range("E51").Value = "  -0.24%  "
